$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 484; this shifts existing rows 484-539 down to 485-540
$ws.Rows.Item(484).Insert()

# Fill the new row 484 with data (same constant columns as surrounding rows, new D/J/K/L/M/O/P)
$ws.Cells.Item(484, 1).Value = 3
$ws.Cells.Item(484, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(484, 3).Value = "Coquimbo"
$ws.Cells.Item(484, 4).Value = 44946
$ws.Cells.Item(484, 5).Value = 5
$ws.Cells.Item(484, 6).Value = 100112017
$ws.Cells.Item(484, 7).Value = "Apio"
$ws.Cells.Item(484, 8).Value = "Americana (o)"
$ws.Cells.Item(484, 9).Value = "Primera"
$ws.Cells.Item(484, 10).Value = 170
$ws.Cells.Item(484, 11).Value = 10000
$ws.Cells.Item(484, 12).Value = 10500
$ws.Cells.Item(484, 13).Value = 10147
$ws.Cells.Item(484, 14).Value = "`$/docena de matas"
$ws.Cells.Item(484, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(484, 16).Value = 1691
$ws.Cells.Item(484, 17).Value = 6
$ws.Cells.Item(484, 18).Value = "Hortaliza"
